$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Automated Document Processing System" "Automated Academic Manuscript Formatter"
Replace-Text "This is the abstract-like section but without the keyword." "Rohit Kumar, Ananya Sharma, Vikram Patel"
Replace-Text "1 System Overview" "Department of Computer Science"
Replace-Text "Content for section 1. This system uses advanced heuristics." "XYZ University"
Replace-Text "2 Technical Approach" "Abstract"
Replace-Text "Content for section 2. We utilize a multi-layered pipeline." "This paper presents an automated system for converting poorly formatted academic manuscripts into publication-ready documents. The system reconstructs documents using deterministic rules and predefined templates, eliminating manual formatting effort and reducing submission errors."
Replace-Text "2.1 Design Considerations" "1 IntroductionAcademic publishing requires strict adherence to formatting guidelines imposed by journals and conferences. Researchers often spend significant time manually adjusting fonts, margins, headings, and references. This process is time-consuming, error-prone, and distracts from core research activities."
Replace-Text "Content for section 2.1. Handling edge cases is priority." "2ethodology"
Replace-Text "3 Experimental Evaluation" "The proposed system uses a pipeline-based architecture in which each stage performs a specific responsibility such as document parsing, structure detection, semantic classification, ad formatting. This modular approach ensures scalability, correctness, and maintainability."
Replace-Text "Content for section 3. We tested on various datasets." "Figure 1: college logo"
Replace-Text "4 Observations and Insights" "Table 1: Comparison of Formatting Methods"
Replace-Text "Content for section 4. The system is extremely robust." "REFERENCES"
Replace-Text "5 Closing Remarks" "[1] John Smith, Automated Document Processing, Journal of Artificial Intelligence, 2021"
Replace-Text "Content for section 5. Future work includes Stage 2 implementation." "[2] Alice Brown, Acaic Formatting Tools and Systems, Publishing Technologies, 2020."
